$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 723.3333
$ws.Range("I9").Value = 468
$ws.Range("K9").Value = 468
$ws.Range("M9").Value = -299

$ws.Range("H43").Value = 2899.5
$ws.Range("J43").Value = 2899.5
$ws.Range("L43").Value = 2899.5
$ws.Range("N43").Value = -3037.5

$ws.Range("H74").Value = 12639.615
$ws.Range("I74").Value = 4052.5
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 4052.5
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -3116.5
$ws.Range("N74").Value = -21872

$ws.Range("H77").Value = 12639.615
$ws.Range("I77").Value = 4052.5
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 20262.5
$ws.Range("L77").Value = 100000
$ws.Range("M77").Value = -15582.5
$ws.Range("N77").Value = -109360

$ws.Range("H116").Value = 4999.8
$ws.Range("I116").Value = 4999.5
$ws.Range("K116").Value = 4999.5
$ws.Range("M116").Value = -1557.5

$ws.Range("H123").Value = 72306.664
$ws.Range("J123").Value = 72306.664
$ws.Range("L123").Value = 72306.664
$ws.Range("N123").Value = -82106.664

$ws.Range("H138").Value = 2180.5386
$ws.Range("I138").Value = 1838
$ws.Range("K138").Value = 5514
$ws.Range("M138").Value = -374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3634.9697
$ws.Range("I32").Value = 3634.9697
$ws.Range("K32").Value = 3634.9697
$ws.Range("M32").Value = -3347.9697

$ws.Range("H101").Value = 54777
$ws.Range("J101").Value = 54777
$ws.Range("L101").Value = 54777
$ws.Range("N101").Value = -61267

$ws.Range("H102").Value = 6142.1665
$ws.Range("I102").Value = 4812.3076
$ws.Range("J102").Value = 9599.799999999999
$ws.Range("K102").Value = 4812.3076
$ws.Range("L102").Value = 9599.799999999999
$ws.Range("M102").Value = -3190.3076
$ws.Range("N102").Value = -12843.8

$ws.Range("H122").Value = 1699.125
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15198.429
$ws.Range("I82").Value = 15198.429
$ws.Range("K82").Value = 15198.429
$ws.Range("M82").Value = -14815.429

$ws.Range("H85").Value = 15198.429
$ws.Range("I85").Value = 15198.429
$ws.Range("K85").Value = 15198.429
$ws.Range("M85").Value = -13872.429

$ws.Range("H99").Value = 4905.875
$ws.Range("I99").Value = 4851.8
$ws.Range("K99").Value = 4851.8
$ws.Range("M99").Value = -3353.8

$ws.Range("H103").Value = 39250
$ws.Range("J103").Value = 39250
$ws.Range("L103").Value = 39250
$ws.Range("N103").Value = -41594

$ws.Range("H107").Value = 6137.154
$ws.Range("I107").Value = 1630.5
$ws.Range("K107").Value = 1630.5
$ws.Range("M107").Value = 289.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1506
$ws.Range("I16").Value = 1419.8
$ws.Range("K16").Value = 1419.8
$ws.Range("M16").Value = -1132.8

$ws.Range("H31").Value = 6376
$ws.Range("J31").Value = 7988.8
$ws.Range("L31").Value = 7988.8
$ws.Range("N31").Value = -8578.799999999999

$ws.Range("H34").Value = 6376
$ws.Range("J34").Value = 7988.8
$ws.Range("L34").Value = 7988.8
$ws.Range("N34").Value = -8392.799999999999

$ws.Range("H55").Value = 18000
$ws.Range("I55").Value = 18000
$ws.Range("K55").Value = 18000
$ws.Range("M55").Value = -17685

$ws.Range("H68").Value = 39858.145
$ws.Range("I68").Value = 13666.667
$ws.Range("J68").Value = 59501.75
$ws.Range("K68").Value = 13666.667
$ws.Range("L68").Value = 59501.75
$ws.Range("M68").Value = -12917.667
$ws.Range("N68").Value = -60999.75

$ws.Range("H71").Value = 39858.145
$ws.Range("I71").Value = 13666.667
$ws.Range("J71").Value = 59501.75
$ws.Range("K71").Value = 41000.001
$ws.Range("L71").Value = 178505.25
$ws.Range("M71").Value = -37256.001
$ws.Range("N71").Value = -185993.25

$ws.Range("H74").Value = 58131
$ws.Range("J74").Value = 62110.668
$ws.Range("L74").Value = 62110.668
$ws.Range("N74").Value = -63858.668

$ws.Range("H77").Value = 58131
$ws.Range("J77").Value = 62110.668
$ws.Range("L77").Value = 186332.004
$ws.Range("N77").Value = -195068.004

$ws.Range("H94").Value = 3147.9285
$ws.Range("J94").Value = 5147.143
$ws.Range("L94").Value = 5147.143
$ws.Range("N94").Value = -6049.143

$ws.Range("H99").Value = 3777.1875
$ws.Range("I99").Value = 3572
$ws.Range("K99").Value = 3572
$ws.Range("M99").Value = -2074

$ws.Range("H107").Value = 347.75
$ws.Range("I107").Value = 412.33334
$ws.Range("J107").Value = 264.7143
$ws.Range("K107").Value = 412.33334
$ws.Range("L107").Value = 264.7143
$ws.Range("M107").Value = 1507.66666
$ws.Range("N107").Value = -4104.7143

$ws.Range("H113").Value = 1506
$ws.Range("I113").Value = 1419.8
$ws.Range("K113").Value = 1419.8
$ws.Range("M113").Value = 750.2

$ws.Range("H126").Value = 3777.1875
$ws.Range("I126").Value = 3572
$ws.Range("K126").Value = 10716
$ws.Range("M126").Value = -8246

$ws.Range("H134").Value = 2129.6924
$ws.Range("I134").Value = 2223.8333
$ws.Range("K134").Value = 6671.499899999999
$ws.Range("M134").Value = -4136.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 147332.95
$ws.Range("I4").Value = 1000000.3
$ws.Range("K4").Value = 3000000.9
$ws.Range("M4").Value = -2999888.9

$ws.Range("H132").Value = 2169.5715
$ws.Range("I132").Value = 1715.3334
$ws.Range("K132").Value = 15438.0006
$ws.Range("M132").Value = -12908.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 19029.5
$ws.Range("I47").Value = 16030
$ws.Range("J47").Value = 22029
$ws.Range("K47").Value = 16030
$ws.Range("L47").Value = 22029
$ws.Range("M47").Value = -15462
$ws.Range("N47").Value = -23165

$ws.Range("H55").Value = 5674.143
$ws.Range("J55").Value = 4444.5
$ws.Range("L55").Value = 4444.5
$ws.Range("N55").Value = -5098.5

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -63744

$ws.Range("H102").Value = 1300.1613
$ws.Range("I102").Value = 904.2143
$ws.Range("K102").Value = 904.2143
$ws.Range("M102").Value = 717.7857

$ws.Range("H122").Value = 153426.64
$ws.Range("I122").Value = 193592.27
$ws.Range("J122").Value = 4240
$ws.Range("K122").Value = 580776.8099999999
$ws.Range("L122").Value = 12720
$ws.Range("M122").Value = -578326.8099999999
$ws.Range("N122").Value = -17620

$ws.Range("H126").Value = 4660.727
$ws.Range("I126").Value = 4270
$ws.Range("K126").Value = 12810
$ws.Range("M126").Value = -10340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 7254.3335
$ws.Range("I32").Value = 8305.200000000001
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 8305.200000000001
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -7988.200000000001
$ws.Range("N32").Value = -2634

$ws.Range("H40").Value = 6319.8096
$ws.Range("I40").Value = 5779.615
$ws.Range("J40").Value = 7197.625
$ws.Range("K40").Value = 5779.615
$ws.Range("L40").Value = 7197.625
$ws.Range("M40").Value = -5643.615
$ws.Range("N40").Value = -7469.625

$ws.Range("H64").Value = 14931.333
$ws.Range("J64").Value = 14931.333
$ws.Range("L64").Value = 14931.333
$ws.Range("N64").Value = -15381.333

$ws.Range("H67").Value = 14931.333
$ws.Range("J67").Value = 14931.333
$ws.Range("L67").Value = 14931.333
$ws.Range("N67").Value = -16491.333

$ws.Range("H93").Value = 1024.5
$ws.Range("I93").Value = 1024.5
$ws.Range("K93").Value = 1024.5
$ws.Range("M93").Value = 223.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 23999.6
$ws.Range("J63").Value = 23999.6
$ws.Range("L63").Value = 23999.6
$ws.Range("N63").Value = -25247.6

$ws.Range("H66").Value = 23999.6
$ws.Range("J66").Value = 23999.6
$ws.Range("L66").Value = 71998.79999999999
$ws.Range("N66").Value = -78238.79999999999

$ws.Range("H101").Value = 6400.6665
$ws.Range("J101").Value = 6400.6665
$ws.Range("L101").Value = 6400.6665
$ws.Range("N101").Value = -12890.6665

$ws.Range("H113").Value = 1165.1666
$ws.Range("I113").Value = 2650
$ws.Range("K113").Value = 7950
$ws.Range("M113").Value = -5780

$ws.Range("H117").Value = 78136.336
$ws.Range("J117").Value = 78136.336
$ws.Range("L117").Value = 78136.336
$ws.Range("N117").Value = -87314.336

$ws.Range("H127").Value = 220000
$ws.Range("I127").Value = 220000
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 220000
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -215040
$ws.Range("N127").ClearContents()
